$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the timing values in the "AQ32 + LA Timing" table (rows 11-17, column J)
$ws.Range("J11").Value = 233
$ws.Range("J12").Value = 69
$ws.Range("J14").Value = 4.9000000000000004
$ws.Range("J17").Value = 1

# Update the selected cell / view state
$ws.Activate()
$ws.Range("J18").Select()
